{"js": "// Sprint 5 Acta: fix the meeting dates (Feb -> Mar) for the 1st and 2nd\n// reunions, and relocate the \"_GoBack\" bookmark from the end of the\n// \"Se entrega las Opciones.\" bullet to just before \"/2020\" in the\n// \"3\u00aa Reuni\u00f3n 01/03/2020\" heading (mirrors the author's manual re-edit\n// of the document, which nudges Word's last-cursor-position bookmark).\n\nconst body = context.document.body;\n\n// --- 1) \"1\u00aa Reuni\u00f3n 25/02/2020\" -> \"1\u00aa Reuni\u00f3n 25/03/2020\" -----------\nlet firstDate = body.search(\"25\\t/02/2020\", { matchCase: false });\nawait context.sync();\nif (firstDate.items.length > 0) {\n  let monthPart = firstDate.items[0].search(\"/02/\", { matchCase: false });\n  await context.sync();\n  let monthDigit = monthPart.items[0].search(\"2\", { matchCase: false });\n  await context.sync();\n  monthDigit.items[0].insertText(\"3\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- 2) \"2\u00aa Reuni\u00f3n 28/02/2020\" -> \"2\u00aa Reuni\u00f3n 28/03/2020\" -----------\nlet secondDate = body.search(\"28\\t/02/2020\", { matchCase: false });\nawait context.sync();\nif (secondDate.items.length > 0) {\n  let monthPart2 = secondDate.items[0].search(\"/02/\", { matchCase: false });\n  await context.sync();\n  let monthDigit2 = monthPart2.items[0].search(\"2\", { matchCase: false });\n  await context.sync();\n  monthDigit2.items[0].insertText(\"3\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- 3) Remove the \"_GoBack\" bookmark from its old location ----------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 4) Re-insert \"_GoBack\" right before \"/2020\" in the 3rd reunion --\nlet thirdDate = body.search(\"01\\t/03/2020\", { matchCase: false });\nawait context.sync();\nif (thirdDate.items.length > 0) {\n  let tailPart = thirdDate.items[0].search(\"/2020\", { matchCase: false });\n  await context.sync();\n  let insertionPoint = tailPart.items[0].getRange(\"Start\");\n  insertionPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Sprint 5 Acta: fix the meeting dates (Feb -> Mar) for the 1st and 2nd\n# reunions, and relocate the \"_GoBack\" bookmark from the end of the\n# \"Se entrega las Opciones.\" bullet to just before \"/2020\" in the\n# \"3\u00aa Reuni\u00f3n 01/03/2020\" heading (mirrors the author's manual re-edit\n# of the document, which nudges Word's last-cursor-position bookmark).\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphContaining($doc, $needle) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        $para = $doc.Paragraphs($i)\n        if ($para.Range.Text -like \"*$needle*\") {\n            return $para.Range\n        }\n    }\n    return $null\n}\n\n# --- 1) \"1\u00aa Reuni\u00f3n 25/02/2020\" -> \"1\u00aa Reuni\u00f3n 25/03/2020\" -----------\n$r1 = Find-ParagraphContaining $d \"Reuni\"\n$p1 = $d.Paragraphs(3).Range\nif ($p1.Text -notlike \"*1*Reuni*\") {\n    $p1 = Find-ParagraphContaining $d \"1\"\n}\n$find = $p1.Find\n$find.ClearFormatting()\n$find.Text = \"/02/\"\nif ($find.Execute()) {\n    $monthDigit = $d.Range($p1.Start + 2, $p1.Start + 3)\n    $monthDigit.Text = \"3\"\n}\n\n# --- 2) \"2\u00aa Reuni\u00f3n 28/02/2020\" -> \"2\u00aa Reuni\u00f3n 28/03/2020\" -----------\n$p2 = $d.Paragraphs(7).Range\n$find2 = $p2.Find\n$find2.ClearFormatting()\n$find2.Text = \"/02/\"\nif ($find2.Execute()) {\n    $monthDigit2 = $d.Range($p2.Start + 2, $p2.Start + 3)\n    $monthDigit2.Text = \"3\"\n}\n\n# --- 3) Remove the \"_GoBack\" bookmark from its old location ----------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- 4) Re-insert \"_GoBack\" right before \"/2020\" in the 3rd reunion --\n$p3 = $d.Paragraphs(11).Range\n$find3 = $p3.Find\n$find3.ClearFormatting()\n$find3.Text = \"/2020\"\nif ($find3.Execute()) {\n    $insertionPoint = $d.Range($p3.Start, $p3.Start)\n    $d.Bookmarks.Add(\"_GoBack\", $insertionPoint)\n}\n"}
